# Reviewed exercises series 1: move "CONTINUE HERE" marker slide from
# position 19 to the end of this block (position 23), and flesh out
# several slides that were still placeholders/stubs.

$p = $ppt.ActivePresentation

# --- 1. Reorder slides: move slide 19 ("CONTINUE HERE") to position 23 ---
$p.Slides.Item(19).MoveTo(23)

# --- 2. Slide now at position 19: "Running WireMock standalone" ---
$s19 = $p.Slides.Item(19)
$body19 = $s19.Shapes.Item(2).TextFrame.TextRange
$body19.Text = "Start WireMock server`r" + `
    "Options: port, keystore, ...`r" + `
    "`r" + `
    "Make mocks permanently available`r" + `
    "For example for multiple teams`r" + `
    "`r" + `
    "Reconfigure mocks via JSON`r" + `
    "`r" + `
    "java -jar wiremock-standalone-2.18.0.jar --port 9876`r" + `
    "`r" + `
    ""

$body19.Paragraphs(2,1).IndentLevel = 2
$body19.Paragraphs(5,1).IndentLevel = 2

$p8 = $body19.Paragraphs(8,1)
$p8.ParagraphFormat.Bullet.Type = 0

$p9 = $body19.Paragraphs(9,1)
$p9.ParagraphFormat.Bullet.Type = 0
$p9.Font.Size = 24
$p9.Font.Italic = $true

$s19.Shapes.Item(2).TextFrame.AutoSize = 2

# --- 3. Slide now at position 20: "Starting and stopping WireMock..." ---
$s20 = $p.Slides.Item(20)
$body20 = $s20.Shapes.Item(2).TextFrame.TextRange
$body20.Text = "Integration in test execution`r" + `
    "`r" + `
    "Mocks in version control (Git, etc.)`r" + `
    "`r" + `
    "JUnit integration using @Rule annotation`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    "Can be used without having to use JUnit as well"

# --- 4. Slide now at position 22: "Exercise time!" (was "Get your hands dirty!") ---
$s22 = $p.Slides.Item(22)
$s22.Shapes.Item(1).TextFrame.TextRange.Text = "Exercise time!"

$body22 = $s22.Shapes.Item(2).TextFrame.TextRange
$body22.Text = "WireMockExercises1.java`r" + `
    "`r" + `
    "Create a number of simple mocks`r" + `
    "`r" + `
    "Exercises are defined in the comments`r" + `
    "`r" + `
    "Verify your solution by running the tests`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    ""
